$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Sema3b", "Nrp1", "ECs", 2, 0.6666666666666666, 2.025829666666667, 6.077489, 0.1667292296771927, 0.1667292296771927, 3, 1, 102.8289443333334, 308.486833, 0.5559120396302444, 0.5559120396302443, 208.3139260224819, 1874.825334202337, 0.09268678613582766, 0.09268678613582765),
    @("ECs", "Sema3b", "Nrp1", "FAPs", 2, 0.6666666666666666, 2.025829666666667, 6.077489, 0.1667292296771927, 0.1667292296771927, 3, 1, 63.66262833333334, 190.987885, 0.3441717873742006, 0.3441717873742006, 128.9696411356406, 1160.726770220765, 0.05738349698552302, 0.05738349698552302),
    @("ECs", "Sema3b", "Nrp1", "sCs", 2, 0.6666666666666666, 2.025829666666667, 6.077489, 0.1667292296771927, 0.1667292296771927, 3, 1, 18.481835, 55.445505, 0.09991617299555507, 0.09991617299555505, 37.44104963743833, 336.969446736945, 0.01665894655584202, 0.01665894655584202),
    @("FAPs", "Sema3b", "Nrp1", "ECs", 3, 1, 5.787456, 17.362368, 0.4763174794741613, 0.4763174794741613, 3, 1, 102.8289443333334, 308.486833, 0.5559120396302444, 0.5559120396302443, 595.1179908556161, 5356.061917700545, 0.2647906215260181, 0.264790621526018),
    @("FAPs", "Sema3b", "Nrp1", "FAPs", 3, 1, 5.787456, 17.362368, 0.4763174794741613, 0.4763174794741613, 3, 1, 63.66262833333334, 190.987885, 0.3441717873742006, 0.3441717873742006, 368.44466032352, 3316.00194291168, 0.1639350382681962, 0.1639350382681962),
    @("FAPs", "Sema3b", "Nrp1", "sCs", 3, 1, 5.787456, 17.362368, 0.4763174794741613, 0.4763174794741613, 3, 1, 18.481835, 55.445505, 0.09991617299555507, 0.09991617299555505, 106.96280686176, 962.66526175584, 0.04759181967994704, 0.04759181967994704),
    @("sCs", "Sema3b", "Nrp1", "ECs", 3, 1, 4.337131333333333, 13.011394, 0.356953290848646, 0.356953290848646, 3, 1, 102.8289443333334, 308.486833, 0.5559120396302444, 0.5559120396302443, 445.9826364416891, 4013.843727975202, 0.1984346319683987, 0.1984346319683986),
    @("sCs", "Sema3b", "Nrp1", "FAPs", 3, 1, 4.337131333333333, 13.011394, 0.356953290848646, 0.356953290848646, 3, 1, 63.66262833333334, 190.987885, 0.3441717873742006, 0.3441717873742006, 276.1131801068544, 2485.01862096169, 0.1228532521204814, 0.1228532521204814),
    @("sCs", "Sema3b", "Nrp1", "sCs", 3, 1, 4.337131333333333, 13.011394, 0.356953290848646, 0.356953290848646, 3, 1, 18.481835, 55.445505, 0.09991617299555507, 0.09991617299555505, 80.15814567599666, 721.4233110839699, 0.035665406759766, 0.035665406759766),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowData[$j]
    }
}

